$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number formats already registered in this workbook's styles (must match
# verbatim - using single-quoted literals so PowerShell does no
# interpolation/escaping - so Excel reuses the existing numFmt entries
# instead of minting new, duplicate ones).
$dateFmt = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'
$moneyFmt = '0.00_);[Red]\(0.00\)'
$countFmt = '0_);[Red]\(0\)'

# Row 28 - 四方坪站, 2026-01-14 (serial 46036)
$ws.Range("A28").Value = 46036
$ws.Range("A28").NumberFormat = $dateFmt
$ws.Range("B28").Value = "四方坪站"
$ws.Range("C28").Value = 14452.04
$ws.Range("C28").NumberFormat = $moneyFmt
$ws.Range("D28").Value = 9273.29
$ws.Range("D28").NumberFormat = $moneyFmt
$ws.Range("E28").Value = 3937.61
$ws.Range("E28").NumberFormat = $moneyFmt
$ws.Range("F28").Value = 648
$ws.Range("F28").NumberFormat = $countFmt

# Row 29 - 高岭站, 2026-01-14 (serial 46036)
$ws.Range("A29").Value = 46036
$ws.Range("A29").NumberFormat = $dateFmt
$ws.Range("B29").Value = "高岭站"
$ws.Range("C29").Value = 4023.47
$ws.Range("C29").NumberFormat = $moneyFmt
$ws.Range("D29").Value = 3410.94
$ws.Range("D29").NumberFormat = $moneyFmt
$ws.Range("E29").Value = 1048.35
$ws.Range("E29").NumberFormat = $moneyFmt
$ws.Range("F29").Value = 151
$ws.Range("F29").NumberFormat = $countFmt

# Move the visible window/selection to where the new rows were entered.
$ws.Range("H28").Select()
